$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.357.17"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.07%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.391.90"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.50%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.72"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.00%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +8.15%  "
$ws.Range("E10").Value = "  +1.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "48.82"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.03%  "
$ws.Range("E12").Value = "  +3.71%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "683.98"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.64"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.939.64"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.49%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "69.431.03"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.23%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.396.05"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.77%  "
$ws.Range("E18").Value = "  +1.76%  "
$ws.Range("E19").Value = "  +1.65%  "
$ws.Range("E20").Value = "  +2.42%  "
$ws.Range("E21").Value = "  +0.55%  "
$ws.Range("E22").Value = "  +1.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.15"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "103.97"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.94"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.10%  "
$ws.Range("E26").Value = "  +1.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.64"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "34.19"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.76"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.98"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.74%  "
$ws.Range("E31").Value = "  +1.75%  "
$ws.Range("E32").Value = "  +10.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "554.32"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.34%  "
$ws.Range("E34").Value = "  +0.65%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "58.12"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.43%  "
$ws.Range("E36").Value = "  +0.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.705.08"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.50%  "
$ws.Range("E38").Value = "  +6.77%  "
$ws.Range("E39").Value = "  +2.34%  "
$ws.Range("E40").Value = "  +1.67%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0₃0708"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.59%  "
$ws.Range("E42").Value = "  +0.96%  "
$ws.Range("E43").Value = "  +0.92%  "
$ws.Range("E44").Value = "  +4.19%  "
$ws.Range("E45").Value = "  -2.84%  "
$ws.Range("E46").Value = "  -0.17%  "
$ws.Range("E47").Value = "  +0.64%  "
$ws.Range("E48").Value = "  +4.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.32"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.98%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.59"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.66%  "
